$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.463.78"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.574.02"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.23"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3399"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07553"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.138"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.35"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.986"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.934"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "1.565.17"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001121"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.97"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06731"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.258"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.39"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.14"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "22.455.99"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.336"
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.579"
$ws.Range("E26").Value = "  -6.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.16"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.87"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.014"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.79"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").Value = "1.741.72"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("E32").Value = "  +4.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.122"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.980"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.829"
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08404"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.379"
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02460"
$ws.Range("E38").Value = "  -3.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2292"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06533"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.452"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.30"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6258"
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.809"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5822"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.086"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.28"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("E50").Value = "  -5.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07321"
$ws.Range("E51").Value = "  -0.07%  "
